$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.087.16'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '3.563.55'
$ws.Range("E3").Value = '  +4.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.05'
$ws.Range("E5").Value = '  +3.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.11'
$ws.Range("E6").Value = '  +3.43%  '
$ws.Range("D7").Value = '3.560.68'
$ws.Range("E7").Value = '  +4.47%  '
$ws.Range("E8").Value = '  +0.29%  '
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  +3.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.00'
$ws.Range("E11").Value = '  -0.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").Value = '  +4.49%  '
$ws.Range("D13").Value = '4.178.71'
$ws.Range("E13").Value = '  +4.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000183'
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.41'
$ws.Range("E15").Value = '  +5.62%  '
$ws.Range("D16").Value = '3.583.71'
$ws.Range("E16").Value = '  +4.20%  '
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").Value = '65.012.55'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.18'
$ws.Range("E19").Value = '  +8.41%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.41'
$ws.Range("E20").Value = '  +7.79%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.87'
$ws.Range("E21").Value = '  +3.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '391.61'
$ws.Range("E22").Value = '  +3.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.580'
$ws.Range("E23").Value = '  +7.96%  '
$ws.Range("D24").Value = '3.716.73'
$ws.Range("E24").Value = '  +4.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.20'
$ws.Range("E25").Value = '  +3.84%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  +13.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.72'
$ws.Range("E28").Value = '  +8.78%  '
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +6.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.30'
$ws.Range("E31").Value = '  +5.08%  '
$ws.Range("D32").Value = '3.582.83'
$ws.Range("E32").Value = '  +4.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.44'
$ws.Range("E33").Value = '  +23.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.03'
$ws.Range("E34").Value = '  +5.40%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.144'
$ws.Range("E36").Value = '  +3.09%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '170.74'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.56'
$ws.Range("E38").Value = '  +8.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.92'
$ws.Range("E39").Value = '  +4.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.04'
$ws.Range("E40").Value = '  +11.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0812'
$ws.Range("E41").Value = '  +8.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.827'
$ws.Range("E42").Value = '  +3.36%  '
$ws.Range("E43").Value = '  +21.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.57'
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.47'
$ws.Range("E46").Value = '  +5.66%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.22'
$ws.Range("E47").Value = '  +10.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.68'
$ws.Range("E48").Value = '  +4.69%  '
$ws.Range("D49").Value = '2.494.60'
$ws.Range("E49").Value = '  +14.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.93'
$ws.Range("E50").Value = '  +8.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '301.63'
$ws.Range("E51").Value = '  +10.08%  '
